# Performance_paper_sheet.xlsx update
# - fill in more "Floyd-Warshall 1k" (column E) datapoints
# - fill in remaining "n/a" datapoints that are now known
# - add a new "Sheet2" with machine/browser CONFIGURATION info

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# --- fill in Floyd-Warshall 1k (column E) datapoints ---
$ws1.Range("E5").Value = 1872
$ws1.Range("E6").Value = 1949
$ws1.Range("E7").Value = 1934
$ws1.Range("E8").Value = 1613
$ws1.Range("E9").Value = 1594
$ws1.Range("E10").Value = 1595
$ws1.Range("E11").Value = 2303

# --- replace remaining "n/a" placeholders with real numbers ---
$ws1.Range("F7").Value = 366
$ws1.Range("F9").Value = 373
$ws1.Range("F11").Value = 471
$ws1.Range("G11").Value = 533

# a couple of columns got manually resized while reviewing the new data
$ws1.Columns.Item(9).ColumnWidth = 15 + 1/6
$ws1.Columns.Item(11).ColumnWidth = 16 + 1/6

# move the active selection, just like a user clicking around after editing
[void]$ws1.Range("A12").Select()

# --- add the new Sheet2 with CONFIGURATION information ---
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Sheet2"

$ws2.Range("B2").Value = "CONFIGURATION"

$ws2.Range("B4").Value = "Machine: Quad core i5-3550, 16GB DDR3 RAM, Sandisk Ultra II 480GB SSD"
$ws2.Range("B5").Value = "OS: Windows 10 Pro"

$ws2.Range("B6").Value = "Firefox:"
$ws2.Range("C6").Value = "53.0.3 (32-bit)"

$ws2.Range("B7").Value = "Chrome: "
$ws2.Range("C7").Value = "58.0.3029.110 (64-bit)"

$ws2.Range("B8").Value = "Edge:"
$ws2.Range("C8").Value = "40.15063.0.0"

$ws2.Range("B9").Value = "NodeJS:"
$ws2.Range("C9").Value = "v7.8.0"

[void]$ws2.Range("B4:C9").Select()

[void]$ws1.Select()
